$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33725
$ws.Range("J3").Value = 33725
$ws.Range("L3").Value = 33725
$ws.Range("N3").Value = -33953

$ws.Range("H5").Value = 135.5
$ws.Range("I5").Value = 161
$ws.Range("J5").Value = 110
$ws.Range("K5").Value = 161
$ws.Range("L5").Value = 110
$ws.Range("M5").Value = -46
$ws.Range("N5").Value = -340

$ws.Range("H11").Value = 773.2143
$ws.Range("I11").Value = 773.2143
$ws.Range("K11").Value = 773.2143
$ws.Range("M11").Value = -633.2143

$ws.Range("H17").Value = 1807.7097
$ws.Range("J17").Value = 1807.7097
$ws.Range("L17").Value = 5423.1291
$ws.Range("N17").Value = -5759.1291

$ws.Range("H40").Value = 1725.1111
$ws.Range("I40").Value = 1545.5555
$ws.Range("J40").Value = 1904.6666
$ws.Range("K40").Value = 1545.5555
$ws.Range("L40").Value = 1904.6666
$ws.Range("M40").Value = -1370.5555
$ws.Range("N40").Value = -2254.6666

$ws.Range("H51").Value = 3096.5557
$ws.Range("I51").Value = 2675.75
$ws.Range("J51").Value = 3169.739
$ws.Range("K51").Value = 2675.75
$ws.Range("L51").Value = 3169.739
$ws.Range("M51").Value = -2191.75
$ws.Range("N51").Value = -4137.739

$ws.Range("H64").Value = 43481212
$ws.Range("I64").Value = 100002000
$ws.Range("J64").Value = 3681.5386
$ws.Range("K64").Value = 100002000
$ws.Range("L64").Value = 3681.5386
$ws.Range("M64").Value = -100001752
$ws.Range("N64").Value = -4177.5386

$ws.Range("H67").Value = 43481212
$ws.Range("I67").Value = 100002000
$ws.Range("J67").Value = 3681.5386
$ws.Range("K67").Value = 100002000
$ws.Range("L67").Value = 3681.5386
$ws.Range("M67").Value = -100001142
$ws.Range("N67").Value = -5397.5386

$ws.Range("H102").Value = 33725
$ws.Range("J102").Value = 33725
$ws.Range("L102").Value = 33725
$ws.Range("N102").Value = -40215

$ws.Range("H138").Value = 1988.6123
$ws.Range("I138").Value = 1512.1111
$ws.Range("J138").Value = 2573.4092
$ws.Range("K138").Value = 4536.3333
$ws.Range("L138").Value = 7720.2276
$ws.Range("M138").Value = 603.6666999999998
$ws.Range("N138").Value = -18000.2276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 168.71428
$ws.Range("I5").Value = 126.2
$ws.Range("J5").Value = 275
$ws.Range("K5").Value = 126.2
$ws.Range("L5").Value = 275
$ws.Range("M5").Value = -14.2
$ws.Range("N5").Value = -499

$ws.Range("H74").Value = 19639.5
$ws.Range("I74").Value = 14760.667
$ws.Range("J74").Value = 34276
$ws.Range("K74").Value = 14760.667
$ws.Range("L74").Value = 34276
$ws.Range("M74").Value = -13886.667
$ws.Range("N74").Value = -36024

$ws.Range("H77").Value = 19639.5
$ws.Range("I77").Value = 14760.667
$ws.Range("J77").Value = 34276
$ws.Range("K77").Value = 73803.33499999999
$ws.Range("L77").Value = 171380
$ws.Range("M77").Value = -69435.33499999999
$ws.Range("N77").Value = -180116

$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 168.71428
$ws.Range("I4").Value = 126.2
$ws.Range("J4").Value = 275
$ws.Range("K4").Value = 126.2
$ws.Range("L4").Value = 275
$ws.Range("M4").Value = -11.2
$ws.Range("N4").Value = -505

$ws.Range("H103").Value = 24885.666
$ws.Range("J103").Value = 24885.666
$ws.Range("L103").Value = 24885.666
$ws.Range("N103").Value = -27229.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 51.923077
$ws.Range("I7").Value = 40.444443
$ws.Range("J7").Value = 77.75
$ws.Range("K7").Value = 40.444443
$ws.Range("L7").Value = 77.75
$ws.Range("M7").Value = 72.55555699999999
$ws.Range("N7").Value = -303.75

$ws.Range("H62").Value = 6047
$ws.Range("I62").Value = 10700.833
$ws.Range("J62").Value = 2944.4443
$ws.Range("K62").Value = 10700.833
$ws.Range("L62").Value = 2944.4443
$ws.Range("M62").Value = -10076.833
$ws.Range("N62").Value = -4192.4443

$ws.Range("H65").Value = 6047
$ws.Range("I65").Value = 10700.833
$ws.Range("J65").Value = 2944.4443
$ws.Range("K65").Value = 53504.165
$ws.Range("L65").Value = 14722.2215
$ws.Range("M65").Value = -50384.165
$ws.Range("N65").Value = -20962.2215

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 466
$ws.Range("I24").Value = 99
$ws.Range("J24").Value = 649.5
$ws.Range("K24").Value = 297
$ws.Range("L24").Value = 1948.5
$ws.Range("M24").Value = -67
$ws.Range("N24").Value = -2408.5

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("N25").Value = 0

$ws.Range("H26").Value = 87.333336
$ws.Range("I26").Value = 71
$ws.Range("J26").Value = 120
$ws.Range("K26").Value = 213
$ws.Range("L26").Value = 360
$ws.Range("M26").Value = 75
$ws.Range("N26").Value = -936

$ws.Range("H29").Value = 27777986
$ws.Range("I29").Value = 110.5
$ws.Range("J29").Value = 33333562
$ws.Range("K29").Value = 331.5
$ws.Range("L29").Value = 100000686
$ws.Range("M29").Value = -54.5
$ws.Range("N29").Value = -100001240

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("N30").Value = 0

$ws.Range("H31").Value = 800
$ws.Range("J31").Value = 800
$ws.Range("L31").Value = 2400
$ws.Range("N31").Value = -2976

$ws.Range("H32").Value = 640.1111
$ws.Range("I32").Value = 90.25
$ws.Range("J32").Value = 1080
$ws.Range("K32").Value = 270.75
$ws.Range("L32").Value = 3240
$ws.Range("M32").Value = 12.25
$ws.Range("N32").Value = -3806

$ws.Range("H33").Value = 219.36363
$ws.Range("I33").Value = 74.5
$ws.Range("J33").Value = 302.14285
$ws.Range("K33").Value = 447
$ws.Range("L33").Value = 1812.8571
$ws.Range("M33").Value = -164
$ws.Range("N33").Value = -2378.8571

$ws.Range("H34").Value = 593.88
$ws.Range("I34").Value = 331.54544
$ws.Range("J34").Value = 800
$ws.Range("K34").Value = 994.63632
$ws.Range("L34").Value = 2400
$ws.Range("M34").Value = -910.63632
$ws.Range("N34").Value = -2568

$ws.Range("H35").Value = 1900
$ws.Range("J35").Value = 1900
$ws.Range("L35").Value = 5700
$ws.Range("N35").Value = -6276

$ws.Range("H36").Value = 66667540
$ws.Range("I36").Value = 1088
$ws.Range("K36").Value = 3264
$ws.Range("M36").Value = -3095

$ws.Range("H39").Value = 3024.5
$ws.Range("J39").Value = 3141.5789
$ws.Range("L39").Value = 9424.736699999999
$ws.Range("N39").Value = -10012.7367

$ws.Range("H109").Value = 353
$ws.Range("J109").Value = 320
$ws.Range("L109").Value = 960
$ws.Range("N109").Value = -3040

$ws.Range("H112").Value = 2614.2856
$ws.Range("I112").Value = 1650
$ws.Range("J112").Value = 3000
$ws.Range("K112").Value = 4950
$ws.Range("L112").Value = 9000
$ws.Range("M112").Value = -3842
$ws.Range("N112").Value = -11216

$ws.Range("H131").Value = 635.5
$ws.Range("I131").Value = 499.85715
$ws.Range("J131").Value = 952
$ws.Range("K131").Value = 1499.57145
$ws.Range("L131").Value = 2856
$ws.Range("M131").Value = 3540.42855
$ws.Range("N131").Value = -12936

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 18599
$ws.Range("J98").Value = 18599
$ws.Range("L98").Value = 18599
$ws.Range("N98").Value = -24589

$ws.Range("H101").Value = 22447.5
$ws.Range("J101").Value = 22447.5
$ws.Range("L101").Value = 22447.5
$ws.Range("N101").Value = -28937.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 34999.5
$ws.Range("J103").Value = 34999.5
$ws.Range("L103").Value = 34999.5
$ws.Range("N103").Value = -37343.5
